$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.12839389494417475
$ws.Range("A2").Value = -0.005999999985487392
$ws.Range("A3").Value = -0.0039999999876538794
$ws.Range("A4").Value = -0.007999999977130301
$ws.Range("A5").Value = -0.0029999999878889128
$ws.Range("A6").Value = 0.015766216628183116
$ws.Range("A7").Value = -0.0099999999685778107
$ws.Range("A8").Value = -0.0099999999684730057
$ws.Range("A9").Value = 0.035945621430697905
$ws.Range("A10").Value = -0.0019999999863085094
$ws.Range("A11").Value = -0.0011870766249826303
$ws.Range("A12").Value = -0.0034999999824019135
$ws.Range("A13").Value = -0.0034999999809279814
$ws.Range("A14").Value = -0.0079999999698605606
$ws.Range("A15").Value = -0.00099999998615984964
$ws.Range("A16").Value = -0.0019999999835360605
$ws.Range("A17").Value = -0.0019999999831430415
$ws.Range("A18").Value = -0.0039999999783146833
$ws.Range("A19").Value = -0.0039999999902109451
$ws.Range("A20").Value = 0.015920646946007366
$ws.Range("A21").Value = -0.0039999999903495009
$ws.Range("A22").Value = -0.0039999999902580186
$ws.Range("A23").Value = -0.0049999999844851928
$ws.Range("A24").Value = -0.019999999947565072
$ws.Range("A25").Value = -0.019999999946863412
$ws.Range("A26").Value = -0.0024999999859822708
$ws.Range("A27").Value = -0.0024999999853037025
$ws.Range("A28").Value = -0.0019999999832762683
$ws.Range("A29").Value = -0.0069999999694703163
$ws.Range("A30").Value = -0.030008428891024952
$ws.Range("A31").Value = -0.0069999999673786562
$ws.Range("A32").Value = -0.0099999999601010359
$ws.Range("A33").Value = -0.003999999974221069
